$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin price/volume data (and row rotation for rows 8-17)
# matching the refreshed symbol list commit.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '315.23'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '3.57%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '35.88'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '0.74%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.140'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '1.39%'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.08103'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '3.07%'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '2.123'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '0.03%'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '8.010'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '1.32%'
$ws.Range('B8').NumberFormat = '@'
$ws.Range('B8').Value = 'GateToken'
$ws.Range('C8').NumberFormat = '@'
$ws.Range('C8').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '4.148'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '1.22%'
$ws.Range('B9').NumberFormat = '@'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').NumberFormat = '@'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.9313'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '1.46%'
$ws.Range('B10').NumberFormat = '@'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').NumberFormat = '@'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1010'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '3.66%'
$ws.Range('B11').NumberFormat = '@'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').NumberFormat = '@'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.1871'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '0.69%'
$ws.Range('B12').NumberFormat = '@'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').NumberFormat = '@'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.09172'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '6.53%'
$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.03616'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '1.64%'
$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.09912'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-0.27%'
$ws.Range('B15').NumberFormat = '@'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').NumberFormat = '@'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.001454'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '0.90%'
$ws.Range('B16').NumberFormat = '@'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').NumberFormat = '@'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.005716'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '1.53%'
$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.469'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '0.05%'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '3.38%'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '-1.63%'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '1.82%'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.147'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '-1.26%'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.2199'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '-0.18%'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.04582'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '0.72%'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '1.33%'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.004703'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '-6.95%'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '-21.71%'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0004518'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '-4.82%'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01973'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '6.64%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.04894'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '3.52%'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.007852'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '4.77%'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1391'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '-0.39%'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.007860'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '1.44%'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '-5.70%'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.01162'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '2.94%'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00006529'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '3.14%'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '0.44%'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '35.41'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '-24.44%'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.001908'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '-4.57%'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00002109'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '0.44%'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0002008'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '0.44%'
